$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.705157251606181
$ws.Cells.Item(2, 3).Value = 0.2235001018007381
$ws.Cells.Item(2, 5).Value = 0.4293450414994595
$ws.Cells.Item(2, 6).Value = 0.4443680307746121
$ws.Cells.Item(2, 7).Value = 0.149740506781761
$ws.Cells.Item(2, 8).Value = 0.3030882085640982
$ws.Cells.Item(2, 9).Value = 0.1919372744917425
$ws.Cells.Item(2, 14).Value = 0.768435978604856
$ws.Cells.Item(2, 15).Value = 0.8039366112629551

$ws.Cells.Item(3, 2).Value = 0.6152365081597111
$ws.Cells.Item(3, 3).Value = 0.2013731579265254
$ws.Cells.Item(3, 5).Value = 0.3745221541837083
$ws.Cells.Item(3, 6).Value = 0.3878228170618172
$ws.Cells.Item(3, 7).Value = 0.1457148568727789
$ws.Cells.Item(3, 8).Value = 0.3048400071515829
$ws.Cells.Item(3, 9).Value = 0.1957187233693567
$ws.Cells.Item(3, 14).Value = 0.7611934265705145
$ws.Cells.Item(3, 15).Value = 0.7988159371562347

$ws.Cells.Item(4, 2).Value = 0.5598181097710722
$ws.Cells.Item(4, 3).Value = 0.1877066222612029
$ws.Cells.Item(4, 5).Value = 0.3409411314559634
$ws.Cells.Item(4, 6).Value = 0.3531389305169483
$ws.Cells.Item(4, 7).Value = 0.1434356675502713
$ws.Cells.Item(4, 8).Value = 0.3061232716694064
$ws.Cells.Item(4, 9).Value = 0.1982350354935178
$ws.Cells.Item(4, 14).Value = 0.7570493373200264
$ws.Cells.Item(4, 15).Value = 0.7965060733484393

$ws.Cells.Item(5, 2).Value = 0.5371841563026578
$ws.Cells.Item(5, 3).Value = 0.1821174582842957
$ws.Cells.Item(5, 5).Value = 0.3272752423178531
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.1425549428299746
$ws.Cells.Item(5, 8).Value = 0.3066983381581778
$ws.Cells.Item(5, 9).Value = 0.1993091804402436
$ws.Cells.Item(5, 14).Value = 0.7554371332630438
$ws.Cells.Item(5, 15).Value = 0.7957736030234202

$ws.Cells.Item(6, 2).Value = 0.5334227984923245
$ws.Cells.Item(6, 3).Value = 0.1811881865137082
$ws.Cells.Item(6, 5).Value = 0.3250071122517255
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.1424115914948132
$ws.Cells.Item(6, 8).Value = 0.3067969731452109
$ws.Cells.Item(6, 9).Value = 0.1994904791842522
$ws.Cells.Item(6, 14).Value = 0.7551740627057626
$ws.Cells.Item(6, 15).Value = 0.7956645608510655

$ws.Cells.Item(7, 2).Value = 0.5595130629048981
$ws.Cells.Item(7, 3).Value = 0.18763132517185
$ws.Cells.Item(7, 5).Value = 0.3407567550554234
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.1434235956612397
$ws.Cells.Item(7, 8).Value = 0.3061308162828311
$ws.Cells.Item(7, 9).Value = 0.1982493247297281
$ws.Cells.Item(7, 14).Value = 0.757027284124149
$ws.Cells.Item(7, 15).Value = 0.796495350778784

$ws.Cells.Item(8, 2).Value = 0.6741963329208431
$ws.Cells.Item(8, 3).Value = 0.2158876150416518
$ws.Cells.Item(8, 5).Value = 0.41042444412561
$ws.Cells.Item(8, 6).Value = 0.4248636149813478
$ws.Cells.Item(8, 7).Value = 0.1483122745402312
$ws.Cells.Item(8, 8).Value = 0.3036490754560361
$ws.Cells.Item(8, 9).Value = 0.1932006454992621
$ws.Cells.Item(8, 14).Value = 0.7658761133396297
$ws.Cells.Item(8, 15).Value = 0.8019972712536401

$ws.Cells.Item(9, 2).Value = 0.8973958482533817
$ws.Cells.Item(9, 3).Value = 0.2706485901978226
$ws.Cells.Item(9, 5).Value = 0.5477585498236408
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.1594439504625882
$ws.Cells.Item(9, 8).Value = 0.3004343427554375
$ws.Cells.Item(9, 9).Value = 0.1848512732778662
$ws.Cells.Item(9, 14).Value = 0.7856165940875002
$ws.Cells.Item(9, 15).Value = 0.8194512608771447

$ws.Cells.Item(10, 2).Value = 1.060290898796268
$ws.Cells.Item(10, 3).Value = 0.3104748056302924
$ws.Cells.Item(10, 5).Value = 0.6492110486444744
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.1685887604123195
$ws.Cells.Item(10, 8).Value = 0.2990858896699251
$ws.Cells.Item(10, 9).Value = 0.1796733173949896
$ws.Cells.Item(10, 14).Value = 0.8015580522005905
$ws.Cells.Item(10, 15).Value = 0.8364033044955761

$ws.Cells.Item(11, 2).Value = 1.13414842104828
$ws.Cells.Item(11, 3).Value = 0.3285024546285911
$ws.Cells.Item(11, 5).Value = 0.6955103280049002
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.1729641330302343
$ws.Cells.Item(11, 8).Value = 0.2986938933893555
$ws.Cells.Item(11, 9).Value = 0.177527695274577
$ws.Cells.Item(11, 14).Value = 0.8091190438736646
$ws.Cells.Item(11, 15).Value = 0.8450260095305282

$ws.Cells.Item(12, 2).Value = 1.162079946571396
$ws.Cells.Item(12, 3).Value = 0.3353159330432618
$ws.Cells.Item(12, 5).Value = 0.7130658523356885
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.1746523489593699
$ws.Cells.Item(12, 8).Value = 0.2985774042529528
$ws.Cells.Item(12, 9).Value = 0.1767455745297291
$ws.Cells.Item(12, 14).Value = 0.8120263207356402
$ws.Cells.Item(12, 15).Value = 0.8484233256438358

$ws.Cells.Item(13, 2).Value = 1.156066050984293
$ws.Cells.Item(13, 3).Value = 0.3338491209366623
$ws.Cells.Item(13, 5).Value = 0.7092839061717768
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.1742873606225288
$ws.Cells.Item(13, 8).Value = 0.298601069452161
$ws.Cells.Item(13, 9).Value = 0.176912663847169
$ws.Cells.Item(13, 14).Value = 0.8113982312582237
$ws.Cells.Item(13, 15).Value = 0.8476857613590312

$ws.Cells.Item(14, 2).Value = 1.13644710938496
$ws.Cells.Item(14, 3).Value = 0.3290632700587821
$ws.Cells.Item(14, 5).Value = 0.6969541609390433
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.1731023928914084
$ws.Cells.Item(14, 8).Value = 0.2986836687669694
$ws.Cells.Item(14, 9).Value = 0.1774627397123716
$ws.Cells.Item(14, 14).Value = 0.8093573456381336
$ws.Cells.Item(14, 15).Value = 0.8453028554932871

$ws.Cells.Item(15, 2).Value = 1.124425112634981
$ws.Cells.Item(15, 3).Value = 0.3261300692371378
$ws.Cells.Item(15, 5).Value = 0.6894048810885209
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.1723806611527863
$ws.Cells.Item(15, 8).Value = 0.2987384275483862
$ws.Cells.Item(15, 9).Value = 0.1778036394268838
$ws.Cells.Item(15, 14).Value = 0.8081129759266759
$ws.Cells.Item(15, 15).Value = 0.8438604898062465

$ws.Cells.Item(16, 2).Value = 1.05545908478058
$ws.Cells.Item(16, 3).Value = 0.309294827095016
$ws.Cells.Item(16, 5).Value = 0.6461884060520902
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.1683071862924663
$ws.Cells.Item(16, 8).Value = 0.2991159726122561
$ws.Cells.Item(16, 9).Value = 0.1798177790149857
$ws.Cells.Item(16, 14).Value = 0.8010701133541573
$ws.Cells.Item(16, 15).Value = 0.83585821580877

$ws.Cells.Item(17, 2).Value = 1.01308696651148
$ws.Cells.Item(17, 3).Value = 0.2989437927301992
$ws.Cells.Item(17, 5).Value = 0.6197156212240316
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.1658636757819067
$ws.Cells.Item(17, 8).Value = 0.2994043746570867
$ws.Cells.Item(17, 9).Value = 0.1811072817977504
$ws.Cells.Item(17, 14).Value = 0.7968284652851594
$ws.Cells.Item(17, 15).Value = 0.8311832008174065

$ws.Cells.Item(18, 2).Value = 0.9886927492303244
$ws.Cells.Item(18, 3).Value = 0.2929817474495451
$ws.Cells.Item(18, 5).Value = 0.6045030059497805
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.1644784798791079
$ws.Cells.Item(18, 8).Value = 0.2995910904707841
$ws.Cells.Item(18, 9).Value = 0.1818687102005079
$ws.Cells.Item(18, 14).Value = 0.7944179006128564
$ws.Cells.Item(18, 15).Value = 0.8285799525021673

$ws.Cells.Item(19, 2).Value = 0.9804294026964726
$ws.Cells.Item(19, 3).Value = 0.2909616680811666
$ws.Cells.Item(19, 5).Value = 0.5993545994325586
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.1640129418626088
$ws.Cells.Item(19, 8).Value = 0.2996578842200819
$ws.Cells.Item(19, 9).Value = 0.1821299010510504
$ws.Cells.Item(19, 14).Value = 0.7936067371722118
$ws.Cells.Item(19, 15).Value = 0.8277132234151452

$ws.Cells.Item(20, 2).Value = 1.017599931992493
$ws.Cells.Item(20, 3).Value = 0.3000465502105953
$ws.Cells.Item(20, 5).Value = 0.6225322506753059
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.1661216933638059
$ws.Cells.Item(20, 8).Value = 0.2993715166601731
$ws.Cells.Item(20, 9).Value = 0.1809679674872555
$ws.Cells.Item(20, 14).Value = 0.7972769848769872
$ws.Cells.Item(20, 15).Value = 0.8316719871478426

$ws.Cells.Item(21, 2).Value = 1.142210677196715
$ws.Cells.Item(21, 3).Value = 0.3304693513779
$ws.Cells.Item(21, 5).Value = 0.7005750683284191
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.173449592504042
$ws.Cells.Item(21, 8).Value = 0.2986585393067429
$ws.Cells.Item(21, 9).Value = 0.1773003429816882
$ws.Cells.Item(21, 14).Value = 0.8099556093788465
$ws.Cells.Item(21, 15).Value = 0.8459991794737221

$ws.Cells.Item(22, 2).Value = 1.223436264342979
$ws.Cells.Item(22, 3).Value = 0.3502753124512026
$ws.Cells.Item(22, 5).Value = 0.7517153785289423
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.1784217532966323
$ws.Cells.Item(22, 8).Value = 0.2983788562352601
$ws.Cells.Item(22, 9).Value = 0.1750805091553964
$ws.Cells.Item(22, 14).Value = 0.8184986705214357
$ws.Cells.Item(22, 15).Value = 0.8561332142624565

$ws.Cells.Item(23, 2).Value = 1.180104812585455
$ws.Cells.Item(23, 3).Value = 0.3397116618502594
$ws.Cells.Item(23, 5).Value = 0.7244079293063379
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.1757511483096721
$ws.Cells.Item(23, 8).Value = 0.2985110460576408
$ws.Cells.Item(23, 9).Value = 0.1762489980889548
$ws.Cells.Item(23, 14).Value = 0.8139156941000607
$ws.Cells.Item(23, 15).Value = 0.8506536467455419

$ws.Cells.Item(24, 2).Value = 1.015559725842309
$ws.Cells.Item(24, 3).Value = 0.2995480281004745
$ws.Cells.Item(24, 5).Value = 0.6212588310387588
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.1660049825576664
$ws.Cells.Item(24, 8).Value = 0.2993863066341049
$ws.Cells.Item(24, 9).Value = 0.1810308889652106
$ws.Cells.Item(24, 14).Value = 0.7970741218731945
$ws.Cells.Item(24, 15).Value = 0.8314507437725638

$ws.Cells.Item(25, 2).Value = 0.8372017315912785
$ws.Cells.Item(25, 3).Value = 0.2559049445981714
$ws.Cells.Item(25, 5).Value = 0.5105173743336024
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.1562646028328842
$ws.Cells.Item(25, 8).Value = 0.3011265071809959
$ws.Cells.Item(25, 9).Value = 0.1869428359086864
$ws.Cells.Item(25, 14).Value = 0.7800224720654541
$ws.Cells.Item(25, 15).Value = 0.8140092118242563

Write-Host "applied changes"